# Update gh-pages output data (generated at 456a3b4).
# Only column F ("人气"/popularity-style counter) values change, on all
# four sheets, for the rows listed in the commit's regenerated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 2-30) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 172
$ws.Range("F3").Value  = 388
$ws.Range("F4").Value  = 1084
$ws.Range("F5").Value  = 28
$ws.Range("F6").Value  = 71
$ws.Range("F9").Value  = 303
$ws.Range("F10").Value = 403
$ws.Range("F13").Value = 336
$ws.Range("F15").Value = 328
$ws.Range("F16").Value = 421
$ws.Range("F17").Value = 5403
$ws.Range("F19").Value = 1516
$ws.Range("F21").Value = 4502
$ws.Range("F22").Value = 4502
$ws.Range("F23").Value = 113
$ws.Range("F25").Value = 1435
$ws.Range("F30").Value = 3772

# --- Sheet "演出" (rows 4, 11, 13) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 95
$ws.Range("F11").Value = 2
$ws.Range("F13").Value = 1

# --- Sheet "本地生活" (rows 2, 4) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9351
$ws.Range("F4").Value = 2100

# --- Sheet "全部类型" (rows 2-46, combined list) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 9351
$ws.Range("F4").Value  = 2100
$ws.Range("F5").Value  = 172
$ws.Range("F6").Value  = 388
$ws.Range("F7").Value  = 1085
$ws.Range("F8").Value  = 28
$ws.Range("F9").Value  = 71
$ws.Range("F12").Value = 303
$ws.Range("F13").Value = 403
$ws.Range("F16").Value = 336
$ws.Range("F20").Value = 95
$ws.Range("F22").Value = 328
$ws.Range("F24").Value = 421
$ws.Range("F25").Value = 5403
$ws.Range("F27").Value = 1516
$ws.Range("F32").Value = 4502
$ws.Range("F33").Value = 4502
$ws.Range("F34").Value = 113
$ws.Range("F36").Value = 1435
$ws.Range("F41").Value = 2
$ws.Range("F43").Value = 1
$ws.Range("F46").Value = 3772
